# Sandbox setup: duplicate "Monthly 10 homes" into a new "Monthly 10 homes 2"
# sheet so the data set can be swapped without touching the original sheet
# or its variant name.

$wb = $excel.ActiveWorkbook

$wsSource = $wb.Worksheets.Item("Monthly 10 homes")

# Copy the source sheet right after itself; this duplicates data, styles,
# column widths, etc. exactly.
$wsSource.Copy([System.Reflection.Missing]::Value, $wsSource)

$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Monthly 10 homes 2"

# Replace the copied values with the new "10 homes 2" data set.
$ws3.Range("B2").Value = 40
$ws3.Range("C2").Value = 41
$ws3.Range("D2").Value = 42
$ws3.Range("E2").Value = 43
$ws3.Range("F2").Value = 44
$ws3.Range("G2").Value = 45
$ws3.Range("H2").Value = 46
$ws3.Range("I2").Value = 47
$ws3.Range("J2").Value = 48
$ws3.Range("K2").Value = 49
$ws3.Range("B4").Value = 406.27199999999999
$ws3.Range("C4").Value = 318.98999999999995
$ws3.Range("D4").Value = 675.154
$ws3.Range("E4").Value = 354.26699999999994
$ws3.Range("F4").Value = 558.93599999999981
$ws3.Range("G4").Value = 257.74599999999998
$ws3.Range("H4").Value = 690.98499999999979
$ws3.Range("I4").Value = 527.72899999999993
$ws3.Range("J4").Value = 535.67000000000007
$ws3.Range("K4").Value = 384.78599999999989
$ws3.Range("L4").Value = 4710.5349999999999
$ws3.Range("B5").Value = 391.24900000000014
$ws3.Range("C5").Value = 271.78699999999992
$ws3.Range("D5").Value = 610.51000000000022
$ws3.Range("E5").Value = 331.23499999999996
$ws3.Range("F5").Value = 409.15500000000009
$ws3.Range("G5").Value = 355.06799999999998
$ws3.Range("H5").Value = 692.18099999999981
$ws3.Range("I5").Value = 432.83100000000007
$ws3.Range("J5").Value = 506.02500000000003
$ws3.Range("K5").Value = 396.22999999999996
$ws3.Range("L5").Value = 4396.2709999999997
$ws3.Range("B6").Value = 319.85600000000005
$ws3.Range("C6").Value = 281.54799999999994
$ws3.Range("D6").Value = 537.07299999999998
$ws3.Range("E6").Value = 228.96600000000004
$ws3.Range("F6").Value = 216.82600000000002
$ws3.Range("G6").Value = 296.69999999999993
$ws3.Range("H6").Value = 634.92500000000007
$ws3.Range("I6").Value = 347.27100000000002
$ws3.Range("J6").Value = 314.5089999999999
$ws3.Range("K6").Value = 291.36200000000002
$ws3.Range("L6").Value = 3469.0360000000005
$ws3.Range("B7").Value = 221.238
$ws3.Range("C7").Value = 326.07299999999992
$ws3.Range("D7").Value = 547.10800000000006
$ws3.Range("E7").Value = 229.934
$ws3.Range("F7").Value = 237.36300000000003
$ws3.Range("G7").Value = 275.96299999999997
$ws3.Range("H7").Value = 710.96200000000022
$ws3.Range("I7").Value = 354.30100000000004
$ws3.Range("J7").Value = 326.47699999999998
$ws3.Range("K7").Value = 330.08099999999985
$ws3.Range("L7").Value = 3559.4999999999995
$ws3.Range("B8").Value = 247.71300000000002
$ws3.Range("C8").Value = 400.83600000000001
$ws3.Range("D8").Value = 553.18399999999997
$ws3.Range("E8").Value = 244.37999999999994
$ws3.Range("F8").Value = 193.68400000000005
$ws3.Range("G8").Value = 221.959
$ws3.Range("H8").Value = 660.11999999999978
$ws3.Range("I8").Value = 397.04199999999997
$ws3.Range("J8").Value = 274.30599999999993
$ws3.Range("K8").Value = 377.53799999999995
$ws3.Range("L8").Value = 3570.7619999999997
$ws3.Range("B9").Value = 297.50299999999999
$ws3.Range("C9").Value = 365.10000000000008
$ws3.Range("D9").Value = 653.93000000000006
$ws3.Range("E9").Value = 341.09199999999998
$ws3.Range("F9").Value = 339.72799999999995
$ws3.Range("G9").Value = 236.94299999999998
$ws3.Range("H9").Value = 695.90700000000004
$ws3.Range("I9").Value = 387.185
$ws3.Range("J9").Value = 285.101
$ws3.Range("K9").Value = 517.99900000000014
$ws3.Range("L9").Value = 4120.4880000000003
$ws3.Range("B10").Value = 325.54700000000003
$ws3.Range("C10").Value = 348.14399999999995
$ws3.Range("D10").Value = 633.2969999999998
$ws3.Range("E10").Value = 404.13799999999998
$ws3.Range("F10").Value = 289.24399999999997
$ws3.Range("G10").Value = 258.39699999999999
$ws3.Range("H10").Value = 731.07199999999989
$ws3.Range("I10").Value = 577.52899999999977
$ws3.Range("J10").Value = 283.29200000000003
$ws3.Range("K10").Value = 717.49099999999976
$ws3.Range("L10").Value = 4568.1509999999989
$ws3.Range("B11").Value = 240.02799999999999
$ws3.Range("C11").Value = 278.899
$ws3.Range("D11").Value = 543.77600000000007
$ws3.Range("E11").Value = 241.11599999999996
$ws3.Range("F11").Value = 277.209
$ws3.Range("G11").Value = 301.53699999999998
$ws3.Range("H11").Value = 646.63799999999992
$ws3.Range("I11").Value = 448.09800000000007
$ws3.Range("J11").Value = 226.14300000000006
$ws3.Range("K11").Value = 385.00799999999998
$ws3.Range("L11").Value = 3588.4519999999998
$ws3.Range("B12").Value = 262.91699999999997
$ws3.Range("C12").Value = 333.04499999999996
$ws3.Range("D12").Value = 556.98500000000013
$ws3.Range("E12").Value = 266.43699999999995
$ws3.Range("F12").Value = 264.45899999999995
$ws3.Range("G12").Value = 324.4199999999999
$ws3.Range("H12").Value = 677.14399999999989
$ws3.Range("I12").Value = 472.14100000000008
$ws3.Range("J12").Value = 282.07200000000012
$ws3.Range("K12").Value = 384.86399999999992
$ws3.Range("L12").Value = 3824.4839999999999
$ws3.Range("B13").Value = 229.59900000000002
$ws3.Range("C13").Value = 297.57500000000005
$ws3.Range("D13").Value = 505.92999999999989
$ws3.Range("E13").Value = 252.50499999999997
$ws3.Range("F13").Value = 208.80099999999999
$ws3.Range("G13").Value = 253.31900000000005
$ws3.Range("H13").Value = 644.601
$ws3.Range("I13").Value = 391.11500000000001
$ws3.Range("J13").Value = 269.27500000000003
$ws3.Range("K13").Value = 330.11200000000002
$ws3.Range("L13").Value = 3382.8319999999999
$ws3.Range("B14").Value = 152.49299999999999
$ws3.Range("C14").Value = 294.72900000000004
$ws3.Range("D14").Value = 534.91200000000015
$ws3.Range("E14").Value = 252.84699999999998
$ws3.Range("F14").Value = 581.01700000000005
$ws3.Range("G14").Value = 288.65700000000004
$ws3.Range("H14").Value = 667.58799999999997
$ws3.Range("I14").Value = 421.49
$ws3.Range("J14").Value = 341.07899999999995
$ws3.Range("K14").Value = 316.48500000000007
$ws3.Range("L14").Value = 3851.297
$ws3.Range("B15").Value = 282.04499999999996
$ws3.Range("C15").Value = 213.60999999999999
$ws3.Range("D15").Value = 669.38600000000008
$ws3.Range("E15").Value = 244.38900000000001
$ws3.Range("F15").Value = 537.41100000000006
$ws3.Range("G15").Value = 316.56799999999993
$ws3.Range("H15").Value = 645.38599999999997
$ws3.Range("I15").Value = 496.23400000000004
$ws3.Range("J15").Value = 261.34699999999998
$ws3.Range("K15").Value = 394.03300000000007
$ws3.Range("L15").Value = 4060.4090000000001
$ws3.Range("B16").Value = 3376.46
$ws3.Range("C16").Value = 3730.3359999999998
$ws3.Range("D16").Value = 7021.2450000000008
$ws3.Range("E16").Value = 3391.3059999999996
$ws3.Range("F16").Value = 4113.8329999999996
$ws3.Range("G16").Value = 3387.2769999999996
$ws3.Range("H16").Value = 8097.509
$ws3.Range("I16").Value = 5252.9660000000003
$ws3.Range("J16").Value = 3905.2960000000003
$ws3.Range("K16").Value = 4825.9889999999996
$ws3.Range("L16").Value = 47102.216999999997
$ws3.Range("B18").Value = 64.875
$ws3.Range("C18").Value = 68.061999999999983
$ws3.Range("D18").Value = 222.16899999999998
$ws3.Range("E18").Value = 85.350999999999999
$ws3.Range("F18").Value = 52.38
$ws3.Range("G18").Value = 63.415999999999983
$ws3.Range("H18").Value = 106.05799999999999
$ws3.Range("I18").Value = 80.919999999999987
$ws3.Range("J18").Value = 109.13099999999997
$ws3.Range("K18").Value = 91.373999999999995
$ws3.Range("L18").Value = 943.73599999999988
$ws3.Range("B19").Value = 101.67499999999997
$ws3.Range("C19").Value = 96.578000000000003
$ws3.Range("D19").Value = 370.40100000000001
$ws3.Range("E19").Value = 126.02299999999998
$ws3.Range("F19").Value = 80.582999999999998
$ws3.Range("G19").Value = 111.73299999999999
$ws3.Range("H19").Value = 152.304
$ws3.Range("I19").Value = 113.61600000000001
$ws3.Range("J19").Value = 217.82599999999999
$ws3.Range("K19").Value = 129.13499999999999
$ws3.Range("L19").Value = 1499.874
$ws3.Range("B20").Value = 119.50399999999999
$ws3.Range("C20").Value = 121.44
$ws3.Range("D20").Value = 466.81799999999993
$ws3.Range("E20").Value = 154.94699999999997
$ws3.Range("F20").Value = 107.56699999999998
$ws3.Range("G20").Value = 185.50200000000001
$ws3.Range("H20").Value = 168.24700000000001
$ws3.Range("I20").Value = 135.96100000000001
$ws3.Range("J20").Value = 283.38700000000006
$ws3.Range("K20").Value = 147.71600000000004
$ws3.Range("L20").Value = 1891.0890000000002
$ws3.Range("B21").Value = 128.70500000000001
$ws3.Range("C21").Value = 138.22
$ws3.Range("D21").Value = 506.53899999999999
$ws3.Range("E21").Value = 188.00899999999999
$ws3.Range("F21").Value = 129.08799999999999
$ws3.Range("G21").Value = 249.62000000000003
$ws3.Range("H21").Value = 179.16300000000004
$ws3.Range("I21").Value = 153.74500000000003
$ws3.Range("J21").Value = 326.00900000000001
$ws3.Range("K21").Value = 152.04999999999995
$ws3.Range("L21").Value = 2151.1480000000001
$ws3.Range("B22").Value = 103.43
$ws3.Range("C22").Value = 120.44500000000001
$ws3.Range("D22").Value = 409.55500000000001
$ws3.Range("E22").Value = 178.89100000000002
$ws3.Range("F22").Value = 116.2
$ws3.Range("G22").Value = 228.75600000000003
$ws3.Range("H22").Value = 150.25200000000001
$ws3.Range("I22").Value = 131.77800000000002
$ws3.Range("J22").Value = 300.93700000000007
$ws3.Range("K22").Value = 130.52199999999999
$ws3.Range("L22").Value = 1870.7660000000003
$ws3.Range("B23").Value = 111.006
$ws3.Range("C23").Value = 136.43800000000002
$ws3.Range("D23").Value = 456.90600000000001
$ws3.Range("E23").Value = 201.34900000000002
$ws3.Range("F23").Value = 130.827
$ws3.Range("G23").Value = 251.40600000000001
$ws3.Range("H23").Value = 155.68400000000003
$ws3.Range("I23").Value = 139.11000000000001
$ws3.Range("J23").Value = 319.82500000000005
$ws3.Range("K23").Value = 141.303
$ws3.Range("L23").Value = 2043.8540000000003
$ws3.Range("B24").Value = 112.99000000000002
$ws3.Range("C24").Value = 133.142
$ws3.Range("D24").Value = 432.70800000000003
$ws3.Range("E24").Value = 188.85600000000002
$ws3.Range("F24").Value = 129.11399999999998
$ws3.Range("G24").Value = 245.78100000000006
$ws3.Range("H24").Value = 156.32499999999999
$ws3.Range("I24").Value = 141.29899999999998
$ws3.Range("J24").Value = 313.95600000000007
$ws3.Range("K24").Value = 137.828
$ws3.Range("L24").Value = 1991.999
$ws3.Range("B25").Value = 98.518000000000015
$ws3.Range("C25").Value = 114.45700000000002
$ws3.Range("D25").Value = 411.76100000000002
$ws3.Range("E25").Value = 157.23799999999997
$ws3.Range("F25").Value = 105.85400000000001
$ws3.Range("G25").Value = 199.86899999999997
$ws3.Range("H25").Value = 137.363
$ws3.Range("I25").Value = 124.30600000000001
$ws3.Range("J25").Value = 257.80099999999999
$ws3.Range("K25").Value = 124.499
$ws3.Range("L25").Value = 1731.6660000000002
$ws3.Range("B26").Value = 111.697
$ws3.Range("C26").Value = 120.01899999999999
$ws3.Range("D26").Value = 425.90500000000009
$ws3.Range("E26").Value = 157.94399999999999
$ws3.Range("F26").Value = 110.45
$ws3.Range("G26").Value = 202.65599999999998
$ws3.Range("H26").Value = 153.88999999999999
$ws3.Range("I26").Value = 133.56999999999996
$ws3.Range("J26").Value = 279.60700000000003
$ws3.Range("K26").Value = 135.64100000000002
$ws3.Range("L26").Value = 1831.3790000000001
$ws3.Range("B27").Value = 91.993999999999986
$ws3.Range("C27").Value = 90.7
$ws3.Range("D27").Value = 356.572
$ws3.Range("E27").Value = 118.15700000000002
$ws3.Range("F27").Value = 79.431999999999988
$ws3.Range("G27").Value = 139.21600000000001
$ws3.Range("H27").Value = 134.75900000000001
$ws3.Range("I27").Value = 105.11500000000001
$ws3.Range("J27").Value = 196.13999999999996
$ws3.Range("K27").Value = 122.11399999999999
$ws3.Range("L27").Value = 1434.1990000000001
$ws3.Range("B28").Value = 80.479999999999961
$ws3.Range("C28").Value = 76.384999999999991
$ws3.Range("D28").Value = 263.3490000000001
$ws3.Range("E28").Value = 98.891999999999982
$ws3.Range("F28").Value = 61.866000000000014
$ws3.Range("G28").Value = 115.74799999999998
$ws3.Range("H28").Value = 114.46799999999999
$ws3.Range("I28").Value = 90.445999999999984
$ws3.Range("J28").Value = 146.62799999999999
$ws3.Range("K28").Value = 103.73099999999999
$ws3.Range("L28").Value = 1151.9929999999999
$ws3.Range("B29").Value = 43.441999999999993
$ws3.Range("C29").Value = 44.253
$ws3.Range("D29").Value = 155.75900000000001
$ws3.Range("E29").Value = 61.987999999999985
$ws3.Range("F29").Value = 33.684999999999995
$ws3.Range("G29").Value = 66.321999999999989
$ws3.Range("H29").Value = 69.669000000000011
$ws3.Range("I29").Value = 57.185999999999986
$ws3.Range("J29").Value = 75.953000000000003
$ws3.Range("K29").Value = 60.731999999999985
$ws3.Range("L29").Value = 668.98900000000003
$ws3.Range("B30").Value = 1168.316
$ws3.Range("C30").Value = 1260.1389999999999
$ws3.Range("D30").Value = 4478.442
$ws3.Range("E30").Value = 1717.6450000000002
$ws3.Range("F30").Value = 1137.046
$ws3.Range("G30").Value = 2060.0250000000001
$ws3.Range("H30").Value = 1678.1820000000002
$ws3.Range("I30").Value = 1407.0519999999999
$ws3.Range("J30").Value = 2827.2000000000003
$ws3.Range("K30").Value = 1476.645
$ws3.Range("L30").Value = 19210.692000000003

# Restore the view state on the original "Monthly 10 homes" sheet (no longer
# the active tab) and set the view/selection for the new sheet + workbook.
$ws2 = $wb.Worksheets.Item("Monthly 10 homes")
$ws2.Activate()
$excel.ActiveWindow.TopLeftCell = $ws2.Range("A7")
$ws2.Range("N9").Select()

$ws3.Activate()
$ws3.Range("C34").Select()
